$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range('G2').Value = 'eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range('G5').Value = 'Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Range('G7').Value = 'afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Range('G8').Value = 'Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Range('G15').Value = 'afaf.abdallah@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg'
$ws.Range('G17').Value = 'Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Range('G18').Value = 'Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range('G19').Value = 'eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range('G22').Value = 'Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Range('G24').Value = 'afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Range('G25').Value = 'Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Range('G32').Value = 'afaf.abdallah@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg'
$ws.Range('G34').Value = 'Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Range('G35').Value = 'Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G45').Value = 'backup@backdoor.com, mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, System'
$ws.Range('G48').Value = 'Aya.alaa-eldein@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Range('G51').Value = 'Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range('G52').Value = 'yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
$ws.Range('G62').Value = 'backup@backdoor.com, mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, System'
$ws.Range('G65').Value = 'Aya.alaa-eldein@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Range('G68').Value = 'Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range('G69').Value = 'yasmin.m.senosy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
$ws.Range('G70').Value = 'eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range('G72').Value = 'Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg'
$ws.Range('G76').Value = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range('G77').Value = 'Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range('G83').Value = 'Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Range('G85').Value = 'Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg'
$ws.Range('G86').Value = 'Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range('G87').Value = 'eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
$ws.Range('G92').Value = 'Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range('G98').Value = 'afaf.abdallah@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, nourhanhosni@med.asu.edu.eg, user@user.com'
$ws.Range('G99').Value = 'user@user.com, Walaa.h.ghanima@med.asu.edu.eg'
$ws.Range('G100').Value = 'Aya.alaa-eldein@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg'
$ws.Range('G102').Value = 'Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg'
$ws.Range('G103').Value = 'Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G106').Value = 'Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Range('G109').Value = 'afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Range('G110').Value = 'afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Range('G116').Value = 'enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg'
$ws.Range('G117').Value = 'enas.omran@med.asu.edu.eg, user@user.com, Amr-Saeed@med.asu.edu.eg
, takyas@med.asu.edu.eg'
$ws.Range('G119').Value = 'Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range('G120').Value = 'Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G123').Value = 'Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Range('G126').Value = 'afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Range('G127').Value = 'afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Range('G133').Value = 'enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg'
$ws.Range('G134').Value = 'enas.omran@med.asu.edu.eg, user@user.com, Amr-Saeed@med.asu.edu.eg
, takyas@med.asu.edu.eg'
$ws.Range('G136').Value = 'Remon.Matta@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, marinasorial@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg'
$ws.Range('G137').Value = 'Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range('G143').Value = 'afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg'
$ws.Range('G144').Value = 'Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Range('G149').Value = 'user@user.com, Walaa.h.ghanima@med.asu.edu.eg'
$ws.Range('G150').Value = 'Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Range('G151').Value = 'randa.rabea@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, System, marian.samir@med.asu.edu.eg'
